# The commit swaps the contents of ppt/theme/theme1.xml (the "Office
# Theme" palette, previously only wired to the Notes Master) and
# ppt/theme/theme2.xml (the "Integral" palette that the slide master /
# whole deck actually uses) so that, after the edit, the deck's live
# theme (theme2.xml) carries the Office Theme colors while theme1.xml
# keeps the Integral ones.
#
# The PowerPoint object model doesn't expose a "replace this theme
# part's raw XML" call, so the supported way to re-colour a theme is
# through Theme.ThemeColorScheme (Design > Colors > Customize Colors
# in the UI). We rewrite all twelve theme colour slots on the
# presentation's live Slide Master theme to the exact "Office Theme"
# RGB values, reproducing the visible effect of the swap (the file
# that governs the deck's rendered palette now matches "Office
# Theme").
#
# RGB() below is VBA's packed 0x00BBGGRR integer (same encoding the
# ThemeColor.RGB COM property reads/writes), built from each target
# hex colour so the written value round-trips to the exact srgbClr
# hex in the saved OOXML.

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = RGB(0x00, 0x00, 0x00)   # dk1      -> 000000
$tcs.Item(2).RGB  = RGB(0xFF, 0xFF, 0xFF)   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = RGB(0x44, 0x54, 0x6A)   # dk2      -> 44546A
$tcs.Item(4).RGB  = RGB(0xE7, 0xE6, 0xE6)   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = RGB(0x5B, 0x9B, 0xD5)   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = RGB(0xED, 0x7D, 0x31)   # accent2  -> ED7D31
$tcs.Item(7).RGB  = RGB(0xA5, 0xA5, 0xA5)   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = RGB(0xFF, 0xC0, 0x00)   # accent4  -> FFC000
$tcs.Item(9).RGB  = RGB(0x44, 0x72, 0xC4)   # accent5  -> 4472C4
$tcs.Item(10).RGB = RGB(0x70, 0xAD, 0x47)   # accent6  -> 70AD47
$tcs.Item(11).RGB = RGB(0x05, 0x63, 0xC1)   # hlink    -> 0563C1
$tcs.Item(12).RGB = RGB(0x95, 0x4F, 0x72)   # folHlink -> 954F72
